# Week 9 slides update (course-schedule.xlsx)
#
# - Row 10 ("Week" 9) on the Week-plan sheet gets a new Prepare (F) cell,
#   and the Monday-class (G) / Wednesday-class (H) content is replaced with
#   the new "IgNobel Results Section" activity material. Row height grows
#   to fit the new text.
# - A hyperlink is added on H10 pointing at the Canvas activity page.
# - View state (zoom + selection) is refreshed on every sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week-plan")

# --- Row 10 content -------------------------------------------------

$ws.Rows.Item(10).RowHeight = 68.75

# F10 (new "Prepare" cell for week 9)
$ws.Range("F10").Value = "- [ ] For Wednesday: On Canvas, find the [**Activity: IgNobel Results Sections** page](https://canvas.unl.edu/courses/192983/pages). Download your assigned paper and skim through it to prepare for class."

# G10 (Monday class)
$ws.Range("G10").Value = "- [Countering misconceptions in results](../slides/09-misconceptions.qmd)"

# H10 (Wednesday class) -- mixed formatting: plain text, a blue hyperlink-styled
# run for the bare URL, then plain text again.
$h10Text = "- [IgNobel Results Section Activity](https://canvas.unl.edu/courses/192983/pages) activity (Canvas) and class discussion"
$ws.Range("H10").Value = $h10Text

$urlText = "https://canvas.unl.edu/courses/192983/pages"
$urlStart = $h10Text.IndexOf($urlText) + 1
$urlLen = $urlText.Length
$ws.Range("H10").Characters($urlStart, $urlLen).Font.Color = 16711680

# Hyperlink on H10 pointing at the Canvas activity listing.
$ws.Hyperlinks.Add($ws.Range("H10"), $urlText, "", "", $urlText)

# --- View state -------------------------------------------------------

foreach ($name in @("Week-plan", "due-dates", "SemesterDates", "Sheet3")) {
    $sheet = $wb.Worksheets.Item($name)
    $sheet.Activate()
    $excel.ActiveWindow.Zoom = 60
}

$ws.Activate()
$ws.Range("H10").Select()
